# Weekly fruit/vegetable price update.
# The underlying source rows for this market/product were reshuffled, so the
# data cells (Fecha, Calidad, Volumen, Precio minimo/maximo/ponderado,
# Unidad de comercializacion, Origen, Precio $/Kg, Kg/unidad) for rows
# 2,3,4,5,6,8 need to take on the values that (in the previous version of the
# sheet) belonged to a different row. Row 7 is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (new) <= old Row 4 values
$ws.Range("D2").Value = 44334
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 100
$ws.Range("N2").Value = 11000
$ws.Range("O2").Value = 12000
$ws.Range("P2").Value = 11500
$ws.Range("Q2").Value = "$/caja 12 kilos granel"
$ws.Range("R2").Value = "Región de O'Higgins"
$ws.Range("S2").Value = 11500
$ws.Range("T2").Value = 1

# Row 3 (new) <= old Row 6 values
$ws.Range("D3").Value = 44742
$ws.Range("L3").Value = "Segunda"
$ws.Range("M3").Value = 100
$ws.Range("N3").Value = 14000
$ws.Range("O3").Value = 15000
$ws.Range("P3").Value = 14500
$ws.Range("Q3").Value = "$/caja 18 kilos granel"
$ws.Range("R3").Value = "Región de O'Higgins"
$ws.Range("S3").Value = 806
$ws.Range("T3").Value = 18

# Row 4 (new) <= old Row 3 values
$ws.Range("D4").Value = 44719
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 50
$ws.Range("N4").Value = 14000
$ws.Range("O4").Value = 15000
$ws.Range("P4").Value = 14400
$ws.Range("Q4").Value = "$/caja 18 kilos granel"
$ws.Range("R4").Value = "Región del Maule"
$ws.Range("S4").Value = 800
$ws.Range("T4").Value = 18

# Row 5 (new) <= old Row 8 values
$ws.Range("D5").Value = 44714
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 100
$ws.Range("N5").Value = 14000
$ws.Range("O5").Value = 15000
$ws.Range("P5").Value = 14500
$ws.Range("Q5").Value = "$/caja 18 kilos granel"
$ws.Range("R5").Value = "Región de O'Higgins"
$ws.Range("S5").Value = 806
$ws.Range("T5").Value = 18

# Row 6 (new) <= old Row 5 values
$ws.Range("D6").Value = 44707
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 60
$ws.Range("N6").Value = 12000
$ws.Range("O6").Value = 13000
$ws.Range("P6").Value = 12500
$ws.Range("Q6").Value = "$/caja 12 kilos empedrada"
$ws.Range("R6").Value = "Provincia de Curicó"
$ws.Range("S6").Value = 1042
$ws.Range("T6").Value = 12

# Row 8 (new) <= old Row 2 values
$ws.Range("D8").Value = 44330
$ws.Range("L8").Value = "Primera"
$ws.Range("M8").Value = 100
$ws.Range("N8").Value = 15000
$ws.Range("O8").Value = 16000
$ws.Range("P8").Value = 15500
$ws.Range("Q8").Value = "$/caja 18 kilos granel"
$ws.Range("R8").Value = "Provincia de Curicó"
$ws.Range("S8").Value = 861
$ws.Range("T8").Value = 18
